$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the "pp" header from D1 to F1, and relabel D1 as "pp1" (halved pp) ---
$ws.Range("F1").Value2 = $ws.Range("D1").Value2
$ws.Range("D1").Value2 = "pp1"

# --- Halve every D column value (rows 2..184) ---
$dRange = $ws.Range("D2:D184")
for ($i = 1; $i -le $dRange.Rows.Count; $i++) {
    $cell = $dRange.Cells.Item($i, 1)
    $cell.Value2 = $cell.Value2 / 2
}

# --- Add the new F column: rounded pp (F2 standalone, then filled in chunks so the
#     shared-formula groups land on the same boundaries Excel itself would produce) ---
$ws.Range("F2").Formula = "=ROUND(D2,0)"
$ws.Range("F3:F66").Formula = "=ROUND(D3,0)"
$ws.Range("F67:F130").Formula = "=ROUND(D67,0)"
$ws.Range("F131:F184").Formula = "=ROUND(D131,0)"

# --- Page setup (paper size / orientation) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection moves to H7 ---
$ws.Range("H7").Select()

$wb.Save()
